$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price (D) and volume-change (E) cells, plus a few
# rows where coin identity (B/C) was reordered, per the source diff.
$ws.Range("D2").Value = '28.308.62'
$ws.Range("E2").Value = '  -2.20%  '
$ws.Range("D3").Value = '1.865.05'
$ws.Range("E3").Value = '  -1.93%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = "'319.78"
$ws.Range("E5").Value = '  -1.39%  '
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").Value = "'0.4398"
$ws.Range("E7").Value = '  -4.32%  '
$ws.Range("D8").Value = "'0.3742"
$ws.Range("E8").Value = '  -1.78%  '
$ws.Range("D9").Value = "'0.07568"
$ws.Range("E9").Value = '  -1.76%  '
$ws.Range("D10").Value = "'0.9409"
$ws.Range("E10").Value = '  -3.33%  '
$ws.Range("E11").Value = '  -2.80%  '
$ws.Range("D12").Value = '1.846.93'
$ws.Range("E12").Value = '  -3.05%  '
$ws.Range("D13").Value = "'6.736"
$ws.Range("E13").Value = '  -2.70%  '
$ws.Range("D14").Value = "'5.489"
$ws.Range("E14").Value = '  -2.83%  '
$ws.Range("D15").Value = "'0.06862"
$ws.Range("E15").Value = '  -2.75%  '
$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = '  -0.18%  '
$ws.Range("D17").Value = "'82.20"
$ws.Range("D18").Value = "'0.000009113"
$ws.Range("E18").Value = '  -3.91%  '
$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").Value = "'16.06"
$ws.Range("E20").Value = '  -3.43%  '
$ws.Range("D21").Value = '28.305.81'
$ws.Range("E21").Value = '  -2.13%  '
$ws.Range("D22").Value = "'5.161"
$ws.Range("E22").Value = '  -2.28%  '
$ws.Range("E23").Value = '  -0.85%  '
$ws.Range("D24").Value = '2.085.56'
$ws.Range("E24").Value = '  -2.58%  '
$ws.Range("D25").Value = "'2.036"
$ws.Range("E25").Value = '  -2.94%  '
$ws.Range("D26").Value = "'154.83"
$ws.Range("E26").Value = '  -1.95%  '
$ws.Range("D27").Value = "'18.43"
$ws.Range("E27").Value = '  -3.34%  '
$ws.Range("D28").Value = "'5.377"
$ws.Range("E28").Value = '  -4.26%  '
$ws.Range("D29").Value = "'114.80"
$ws.Range("E29").Value = '  -2.20%  '
$ws.Range("D30").Value = "'1.735"
$ws.Range("E30").Value = '  -5.74%  '
$ws.Range("D31").Value = "'0.09063"
$ws.Range("E31").Value = '  -2.01%  '
$ws.Range("D32").Value = "'0.8078"
$ws.Range("E32").Value = '  -5.79%  '
$ws.Range("D33").Value = "'4.870"
$ws.Range("E33").Value = '  -4.29%  '
$ws.Range("D34").Value = "'1.179"
$ws.Range("D35").Value = "'2.947"
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D36").Value = "'1.002"
$ws.Range("D37").Value = "'1.128"
$ws.Range("E37").Value = '  -1.10%  '
$ws.Range("E38").Value = '  -3.26%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = "'0.01960"
$ws.Range("E39").Value = '  -3.47%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = "'3.007"
$ws.Range("E40").Value = '  +8.99%  '
$ws.Range("D41").Value = "'7.185"
$ws.Range("D42").Value = "'0.5274"
$ws.Range("E42").Value = '  -3.73%  '
$ws.Range("D43").Value = "'0.1681"
$ws.Range("D44").Value = "'8.855"
$ws.Range("E44").Value = '  -4.72%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").Value = "'2.079"
$ws.Range("E45").Value = '  +0.84%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = "'0.06786"
$ws.Range("E46").Value = '  -0.54%  '
$ws.Range("D47").Value = "'0.4904"
$ws.Range("E47").Value = '  -4.87%  '
$ws.Range("D48").Value = "'0.000002546"
$ws.Range("E48").Value = '  -0.69%  '
$ws.Range("D49").Value = "'10.58"
$ws.Range("E49").Value = '  -4.92%  '
$ws.Range("D50").Value = "'107.77"
$ws.Range("E50").Value = '  -2.07%  '
$ws.Range("D51").Value = "'1.693"
$ws.Range("E51").Value = '  -4.19%  '

# Values like "82.20" or "1.004" look numeric to Excel and would be
# auto-converted (dropping meaningful trailing zeros / reformatting),
# so they were entered with a leading apostrophe to force text. Excel
# marks such cells with a "quote prefix" style; reset back to Normal
# so the cells keep the workbook's original (unstyled) appearance.
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
